# Update Sheet1 values for rows 2-25, columns C,D,E,F,G,H,K,L,M
# as described by the commit "case with 380 kV done"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$colLetters = @("C", "D", "E", "F", "G", "H", "K", "L", "M")

$data = @(
    @(2, 0.3653238666359044, 0.08059530152428707, 0.13737175907238, 1.686630687589314, 1.614907548290091, 1.428435217505793, 1.581237214977477, 0.1083823394979504, 0.5138977588375511),
    @(3, 0.3622542047269661, 0.0813245864336043, 0.1372721066200171, 1.682074153880379, 1.612010330419267, 1.434474655793906, 1.446925348307275, 0.1089432774486525, 0.4867074740496449),
    @(4, 0.360560720881864, 0.08179503714347902, 0.1372777038520319, 1.680480148603792, 1.611470734099527, 1.439021226938323, 1.364900502348007, 0.1093339229114818, 0.4702352898092883),
    @(5, 0.3599186713604325, 0.08199245552028955, 0.1372967645639207, 1.680132193302072, 1.611561158678995, 1.441084337214988, 1.331586669326043, 0.1095047319338676, 0.4635788003946146),
    @(6, 0.3598149608813088, 0.0820255815026254, 0.1373009426492651, 1.680092601094188, 1.611594877291878, 1.44143960605534, 1.32606171083637, 0.1095337961462874, 0.4624768827681294),
    @(7, 0.3605518674471995, 0.08179767648695524, 0.1372778929885534, 1.680474236161174, 1.611470698903688, 1.439048199705695, 1.364450766318242, 0.109336179470187, 0.4701452910173529),
    @(8, 0.3642257130900219, 0.08084205805583622, 0.1373235291478387, 1.684809109752095, 1.613650569079837, 1.430343340838178, 1.534834754488372, 0.1085661524548307, 0.504476329334409),
    @(9, 0.372950954620535, 0.07914767006918133, 0.1379437290642009, 1.702912074937302, 1.627820771337383, 1.419946777585892, 1.872470978036972, 0.1074232633736329, 0.5735693463577221),
    @(10, 0.380293943980746, 0.07801190266244262, 0.1387243118438342, 1.722141380791129, 1.644353996551985, 1.416408736975058, 2.122700263217212, 0.1068079381647209, 0.6254208237557464),
    @(11, 0.3838382332105823, 0.0775188339741959, 0.1391502976783912, 1.732193281476071, 1.653224467541889, 1.415696575298199, 2.237014847267176, 0.106576853906418, 0.6492484866796104),
    @(12, 0.3852097697495083, 0.07733551185860676, 0.1393218245764096, 1.736188518448159, 1.656779064793255, 1.415556486424578, 2.28037252678007, 0.1064963801221346, 0.6583060381709913),
    @(13, 0.3849130765537154, 0.07737484269153327, 0.1392844285259436, 1.73531965606665, 1.656004797953983, 1.415580884506369, 2.271031601317873, 0.106513398626273, 0.6563537964919988),
    @(14, 0.3839504809147627, 0.07750368397047325, 0.1391642044490204, 1.73251818111487, 1.653512978835494, 1.415682449681412, 2.240580524223958, 0.1065700922678161, 0.6499929647886233),
    @(15, 0.3833646929006989, 0.07758304473199473, 0.1390918946888249, 1.730826820824973, 1.652012177622908, 1.415761554565279, 2.221937365333815, 0.1066057349727245, 0.6461012689544816),
    @(16, 0.3800664255269055, 0.07804460060300755, 0.1386979010120228, 1.721510818163082, 1.643801572056645, 1.416473382291684, 2.115239243104554, 0.1068240235946654, 0.6238684622343555),
    @(17, 0.3780953226584813, 0.0783337931610486, 0.1384743710943859, 1.716130711682567, 1.63911131399982, 1.417140284721569, 2.049907091228022, 0.1069704509945915, 0.6102908626154147),
    @(18, 0.3769807868433759, 0.07850235020847762, 0.138352474284499, 1.713158918443199, 1.636540525885181, 1.417608294822855, 2.01237534400309, 0.1070592676206772, 0.6025039915511456),
    @(19, 0.3766067181375092, 0.07855980225258818, 0.1383123473644652, 1.712173757067504, 1.635691855064522, 1.417781238181306, 1.999675579591553, 0.1070901283450993, 0.5998713734991128),
    @(20, 0.378303163310477, 0.07830277825147292, 0.1384974755936277, 1.716690725187675, 1.639597453011532, 1.41706055010215, 2.056857095964631, 0.1069543878566535, 0.6117338817794717),
    @(21, 0.3842324203791918, 0.07746574810245477, 0.1391992397664374, 1.733335908389122, 1.65423956758039, 1.415649095977443, 2.249522867666997, 0.1065532490128263, 0.6518603580180127),
    @(22, 0.3882788790444067, 0.07693847341137783, 0.1397174358398843, 1.745315546834391, 1.66494950687408, 1.415482249152745, 2.375844455235324, 0.1063320795110094, 0.678286618106668),
    @(23, 0.3861035071125229, 0.07721808077679881, 0.1394354087919076, 1.738820633891251, 1.659128560919299, 1.415501973206517, 2.308387450811267, 0.1064463667947315, 0.6641640063576233),
    @(24, 0.3782091403803065, 0.07831679293862237, 0.1384870094511435, 1.716437165329907, 1.639377278122055, 1.417096334641741, 2.053714909786834, 0.1069616355690783, 0.6110814334109165),
    @(25, 0.3704271769916119, 0.07958686325399356, 0.1377189909239824, 1.696978996612586, 1.622918714986014, 1.422041583310005, 1.780752751817545, 0.1076930860368641, 0.5546873876687926)
)

foreach ($rowData in $data) {
    $rowNum = $rowData[0]
    for ($i = 0; $i -lt $colLetters.Length; $i++) {
        $col = $colLetters[$i]
        $val = $rowData[$i + 1]
        $ws.Range("$col$rowNum").Value = $val
    }
}
